# Auto update: 2025-12-05 17:31:00
# Apply the refreshed screener values for the 국장_반도체_분석 sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 / Row 7: labels were swapped (DB HiTek <-> 240810.KS combo ticker) ---
$ws.Range("B6").Value = "DB HiTek"
$ws.Range("C6").Value = "000990.KS"
$ws.Range("B7").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C7").Value = "240810.KS"

# --- Row 2 (SamsungElec combo row) ---
$ws.Range("K2").Value = 63.6
$ws.Range("N2").Value = 54.77309453746771

# --- Row 3 (SamsungElec) ---
$ws.Range("D3").Value = 108400
$ws.Range("E3").Value = 61.6
$ws.Range("F3").Value = 7.86
$ws.Range("G3").Value = 50
$ws.Range("K3").Value = 55.4
$ws.Range("N3").Value = 54.77309453746771

# --- Row 4 (SK hynix combo row) ---
$ws.Range("K4").Value = 51.6
$ws.Range("N4").Value = 54.77309453746771

# --- Row 5 (SK hynix) ---
$ws.Range("D5").Value = 544000
$ws.Range("E5").Value = 33.8
$ws.Range("F5").Value = 2.64
$ws.Range("K5").Value = 48.8
$ws.Range("N5").Value = 54.77309453746771

# --- Row 6 (now DB HiTek) ---
$ws.Range("D6").Value = 64800
$ws.Range("E6").Value = 33.9
$ws.Range("F6").Value = 1.89
$ws.Range("H6").Value = 50
$ws.Range("J6").Value = 63
$ws.Range("K6").Value = 40.8
$ws.Range("N6").Value = 54.77309453746771

# --- Row 7 (now 240810.KS combo) ---
$ws.Range("D7").Value = 61000
$ws.Range("E7").Value = 36.7
$ws.Range("F7").Value = 1.16
$ws.Range("H7").Value = 60
$ws.Range("I7").Value = 46
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 40.8
$ws.Range("N7").Value = 54.77309453746771
